# Adds 4 new "BETA" stage Australia users (rows 14-17) to the Users sheet,
# mirroring the existing B/C/D/E/F data-row layout (Stage, Id, Email, 1234567, Country).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row data: Stage, Id, Email, <const>, Country
$rows = @(
    @{ Row = 14; Stage = "BETA"; Id = "70000008066"; Email = "UwLyX+test@RCB.com"; Const = "1234567"; Country = "Australia" },
    @{ Row = 15; Stage = "BETA"; Id = "70000028174"; Email = "yFsOM+test@ZRN.com"; Const = "1234567"; Country = "Australia" },
    @{ Row = 16; Stage = "BETA"; Id = "70000014689"; Email = "bELDl+test@FPP.com"; Const = "1234567"; Country = "Australia" },
    @{ Row = 17; Stage = "BETA"; Id = "70000011098"; Email = "qgsEZ+test@sTT.com"; Const = "1234567"; Country = "Australia" }
)

$firstRow = 14
$lastRow = 17

# Force the numeric-looking Id / Const columns to be stored as text (shared
# strings), matching the source data, instead of being auto-detected as
# numbers. Clearing formats afterwards drops the temporary "@" number format
# again so the cells keep their plain (unstyled) look.
$ws.Range("C$($firstRow):C$($lastRow)").NumberFormat = "@"
$ws.Range("E$($firstRow):E$($lastRow)").NumberFormat = "@"

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 2).Value = $r.Stage
    $ws.Cells.Item($row, 3).Value = $r.Id
    $ws.Cells.Item($row, 4).Value = $r.Email
    $ws.Cells.Item($row, 5).Value = $r.Const
    $ws.Cells.Item($row, 6).Value = $r.Country
}

$ws.Range("C$($firstRow):C$($lastRow)").ClearFormats()
$ws.Range("E$($firstRow):E$($lastRow)").ClearFormats()

Write-Host "Added rows 14-17"
